$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.150250911712646
$ws.Range("B1").Value = 2.581364631652832
$ws.Range("C1").Value = 9.594318389892578
$ws.Range("D1").Value = 2.132050752639771
$ws.Range("E1").Value = 1.246553540229797
